$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All Test Cases")
$ws.Activate()
$ws.Rows.Item(37).Select() | Out-Null
$ws.Rows.Item(37).Delete()
